# Regenerate orders with updated distance/size codes.
# Mapping applied to every text cell in the used range:
#   D80 -> D86
#   D64 -> D69
#   D51 -> D55
#   S30 -> S31
# (S20 / S25 are left unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$firstRow = $used.Row
$firstCol = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($firstRow + $r, $firstCol + $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
